# Updating FoodLog file 22/05/2018
# Adds 8 new days of food-log data (rows 134:141) below the existing table,
# extends the two helper formula columns (I = WaterTargetAchieved,
# J = UnderEaten), widens column B slightly and moves the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows: row, DateSerial, CaloriesIn, Fat, Fiber, Carbs, Sodium, Protein, Water
$newRows = @(
    @(134, 43234, 2041, 67,  31, 169, 3468, 187, 4000),
    @(135, 43235, 1507, 35,  25, 180, 2938, 118, 3500),
    @(136, 43236, 1681, 47,  13, 131, 3242, 169, 3250),
    @(137, 43237, 1522, 43,  20, 143, 2691, 136, 3500),
    @(138, 43238, 1727, 46,  22, 154, 3362, 170, 3250),
    @(139, 43239, 2062, 101, 24, 178, 3590, 117, 1250),
    @(140, 43240, 1861, 45,  21, 218, 3873, 139, 1000),
    @(141, 43241, 1884, 57,  27, 162, 3415, 178, 4000)
)

# Copy the formatting of the last existing row (133) down onto the new rows
# so number formats / fonts / borders match the rest of the table.
$ws.Range("A133:J133").Copy()
$ws.Range("A134:J141").PasteSpecial(-4122)

# Column B on the new rows uses the plain integer "Comma" style seen earlier
# in the sheet (e.g. B56) rather than the special style used on rows 127:133.
$ws.Range("B56").Copy()
$ws.Range("B134:B141").PasteSpecial(-4122)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]   # A - Date
    $ws.Cells.Item($r, 2).Value = $row[2]   # B - Calories In
    $ws.Cells.Item($r, 3).Value = $row[3]   # C - Fat (g)
    $ws.Cells.Item($r, 4).Value = $row[4]   # D - Fiber (g)
    $ws.Cells.Item($r, 5).Value = $row[5]   # E - Carbs (g)
    $ws.Cells.Item($r, 6).Value = $row[6]   # F - Sodium (mg)
    $ws.Cells.Item($r, 7).Value = $row[7]   # G - Protein (g)
    $ws.Cells.Item($r, 8).Value = $row[8]   # H - Water (ml)

    $ws.Cells.Item($r, 9).Formula = "=IF(H$r>=2200,""Yes"",""No"")"
    $ws.Cells.Item($r, 10).Formula = "=IF(B$r<=1800,""Yes"",""No"")"
}

# Column B is slightly wider to fit the new values.
$ws.Columns("B:B").ColumnWidth = 9.3

# Move the active selection like in the saved workbook.
$ws.Range("L137").Select()
